$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

# Header row (row 1) - new headers "PD" (Q1) and "N2" (R1)
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"

# Match the formatting of the existing header cells (bold, centered, bordered)
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-48 for column Q (PD) and R (N2)
$qValues = @{
    2  = "AGU-J"
    3  = "PUE-G"
    4  = "CLI-G"
    5  = "VCR-E"
    6  = "CLI-N"
    7  = "CON-H"
    8  = "RET-C"
    9  = "VCR-N"
    10 = "ALM-B"
    11 = "CLI-N"
    12 = "PCH-G"
    13 = "ATH-P"
    14 = "ATH-P"
    15 = "PUE-B"
    16 = "PPT-F"
    17 = "PPT-K"
    18 = "ALM-G"
    19 = "ALM-F"
    20 = "CON-A"
    21 = "AGU-F"
    22 = "VCR-H"
    23 = "PAV-Q"
    24 = "COG-C"
    25 = "ALM-L"
    26 = "ALM-I"
    27 = "CON-B"
    28 = "CLI-L"
    29 = "PAV-T"
    30 = "CLI-B"
    31 = "ALM-J"
    32 = "ALM-C"
    33 = "PPT-S"
    34 = "ATH-J"
    35 = "NRA-A"
    36 = "AGU-B"
    37 = "PUE-O"
    38 = "ATH-C"
    39 = "PUE-E"
    40 = "PUE-K"
    41 = "PAV-U"
    42 = "VCR-I"
    43 = "AGU-N"
    44 = "CON-E"
    45 = "BLO-G"
    46 = "RET-S"
    47 = "PUE-O"
    48 = "BLO-B"
}

$rValues = @{
    2  = "Fuera de Poligono OVL"
    3  = "Fuera de Poligono OVL"
    4  = "Fuera de Poligono OVL"
    5  = "Fuera de Poligono OVL"
    6  = "Fuera de Poligono OVL"
    7  = "Fuera de Poligono OVL"
    8  = "Fuera de Poligono OVL"
    9  = "Fuera de Poligono OVL"
    10 = "Fuera de Poligono OVL"
    11 = "Fuera de Poligono OVL"
    12 = "Fuera de Poligono OVL"
    13 = "Fuera de Poligono OVL"
    14 = "Fuera de Poligono OVL"
    15 = "Fuera de Poligono OVL"
    16 = "ARATO-25058.PO.1PPT"
    17 = "Fuera de Poligono OVL"
    18 = "Fuera de Poligono OVL"
    19 = "Fuera de Poligono OVL"
    20 = "Fuera de Poligono OVL"
    21 = "Fuera de Poligono OVL"
    22 = "Fuera de Poligono OVL"
    23 = "Fuera de Poligono OVL"
    24 = "Fuera de Poligono OVL"
    25 = "Fuera de Poligono OVL"
    26 = "Fuera de Poligono OVL"
    27 = "Fuera de Poligono OVL"
    28 = "Fuera de Poligono OVL"
    29 = "Fuera de Poligono OVL"
    30 = "Fuera de Poligono OVL"
    31 = "Fuera de Poligono OVL"
    32 = "Fuera de Poligono OVL"
    33 = "Fuera de Poligono OVL"
    34 = "Fuera de Poligono OVL"
    35 = "Fuera de Poligono OVL"
    36 = "Fuera de Poligono OVL"
    37 = "ARATO-25058.PO.1PUE"
    38 = "Fuera de Poligono OVL"
    39 = "Fuera de Poligono OVL"
    40 = "Fuera de Poligono OVL"
    41 = "Fuera de Poligono OVL"
    42 = "Fuera de Poligono OVL"
    43 = "Fuera de Poligono OVL"
    44 = "Fuera de Poligono OVL"
    45 = "Fuera de Poligono OVL"
    46 = "Fuera de Poligono OVL"
    47 = "ARATO-25058.PO.1PUE"
    48 = "Fuera de Poligono OVL"
}

for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 17).Value = $qValues[$row]
    $ws.Cells.Item($row, 18).Value = $rValues[$row]
}
